$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 8.460500000000003
$ws.Range("D5").Value = -7.451100000000005
$ws.Range("A8").Value = -22.55560000000002
$ws.Range("D8").Value = -8.296999999999999
$ws.Range("A10").Value = -22.11129999999999
$ws.Range("B11").Value = 6.9245
$ws.Range("A12").Value = -21.49329999999998
$ws.Range("B12").Value = 5.370099999999997
$ws.Range("D12").Value = -5.852900000000002
$ws.Range("D13").Value = -8.643999999999991
$ws.Range("B15").Value = 5.848299999999996
$ws.Range("D15").Value = -8.103499999999999
$ws.Range("B17").Value = 4.825199999999999
$ws.Range("A18").Value = -22.1659
$ws.Range("D21").Value = -7.556599999999992
$ws.Range("A25").Value = -22.02670000000002
$ws.Range("D25").Value = -7.869799999999999
$ws.Range("B26").Value = 5.109500000000001
$ws.Range("B27").Value = 6.704900000000004
$ws.Range("B28").Value = 6.6051
$ws.Range("B32").Value = 7.7811
$ws.Range("D32").Value = -8.077699999999991
$ws.Range("D36").Value = -6.903300000000006
$ws.Range("A37").Value = -21.79949999999999
$ws.Range("B37").Value = 6.749399999999997
$ws.Range("D38").Value = -7.693099999999998
$ws.Range("B41").Value = 9.014300000000002
$ws.Range("D41").Value = -8.009599999999997
$ws.Range("B47").Value = 6.717300000000003
$ws.Range("D50").Value = -8.013400000000001
$ws.Range("B51").Value = 5.556
$ws.Range("D52").Value = -7.789500000000003
$ws.Range("A55").Value = -21.64299999999999
$ws.Range("D59").Value = -8.272200000000002
$ws.Range("B65").Value = 5.920200000000001
$ws.Range("D67").Value = -7.435499999999995
$ws.Range("A68").Value = -21.4734
$ws.Range("B73").Value = 8.147899999999998
$ws.Range("A77").Value = -19.76359999999999
$ws.Range("A78").Value = -19.79999999999997
$ws.Range("A79").Value = -20.13029999999998
$ws.Range("A80").Value = -19.1915
$ws.Range("A81").Value = -22.24990000000001
$ws.Range("A82").Value = -21.91910000000001
$ws.Range("A84").Value = -21.8596
$ws.Range("B84").Value = 6.352600000000001
$ws.Range("D84").Value = -8.307599999999997
$ws.Range("B85").Value = 6.350499999999998
$ws.Range("D86").Value = -8.183600000000002
$ws.Range("D88").Value = -7.9537
$ws.Range("B89").Value = 4.4004
$ws.Range("D89").Value = -8.096600000000002
$ws.Range("B93").Value = 5.449699999999996
$ws.Range("B95").Value = 5.348200000000005
$ws.Range("D95").Value = -7.726399999999996
$ws.Range("B98").Value = 7.532600000000002
$ws.Range("B99").Value = 6.8111
$ws.Range("A101").Value = -21.8891
$ws.Range("B101").Value = 6.168099999999998
$ws.Range("A102").Value = -21.97870000000002
$ws.Range("B102").Value = 5.479699999999998
$ws.Range("D105").Value = -7.669300000000004
